# LeetCode Practice Tracker update
# - Inserts 4 newly solved problems ("2 Sum", "Container With Most Water",
#   "3 Sum", "3 Sum Closet") at the top of the log (rows 2-5).
# - Pushes the previously existing 17 "To Do" rows down to rows 27-43
#   (leaving a gap, matching the author's original layout).
# - Re-applies the centered alignment style the author used, plus a
#   centered short-date format on the "Date Solved" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Relocate the existing 17 data rows (old rows 2-18) down to rows
#        27-43 by copying values cell-by-cell (bottom-up, so source and
#        destination never collide) and then clearing the vacated cells.
#        Columns A-E are the only ones populated on those rows.
for ($i = 18; $i -ge 2; $i--) {
    $destRow = $i + 25
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $ws.Cells.Item($i, $c).Value2
    }
}
for ($i = 2; $i -le 18; $i++) {
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($i, $c).ClearContents()
    }
}

# --- 2. Write the four newly solved problems into the freed rows 2-5.
#        Columns: A=ID B=Topic C=Problem Name D=Difficulty E=Status
#                 F=Date Solved G=Time Complexity H=Space Complexity
#                 I=Approach Summary

$newRows = @(
    @{ Id = 1; Topic = "Array"; Problem = "2 Sum"; Difficulty = "Easy"; Status = "Done"; Date = 45861; Time = "0(n)"; Space = "O(n)"; Approach = "HashMap" },
    @{ Id = 2; Topic = "Array"; Problem = "Container With Most Water"; Difficulty = "Medium"; Status = "Done"; Date = 45861; Time = "O(n)"; Space = "O(1)"; Approach = "Two Pointers" },
    @{ Id = 3; Topic = "Array"; Problem = "3 Sum"; Difficulty = "Medium"; Status = "Done"; Date = 45863; Time = "O(n ^ 2)"; Space = "O(1)"; Approach = "Two Pointers" },
    @{ Id = 4; Topic = "Array"; Problem = "3 Sum Closet"; Difficulty = "Medium"; Status = "Done"; Date = 45863; Time = "O(n ^ 2)"; Space = "O(1)"; Approach = "Two Pointers" }
)

$r = 2
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row.Id
    $ws.Cells.Item($r, 2).Value = $row.Topic
    $ws.Cells.Item($r, 3).Value = $row.Problem
    $ws.Cells.Item($r, 4).Value = $row.Difficulty
    $ws.Cells.Item($r, 5).Value = $row.Status
    $ws.Cells.Item($r, 6).Value = $row.Date
    $ws.Cells.Item($r, 7).Value = $row.Time
    $ws.Cells.Item($r, 8).Value = $row.Space
    $ws.Cells.Item($r, 9).Value = $row.Approach
    $r = $r + 1
}

# --- 3. Re-apply the centered alignment across the whole used body
#        (columns A-L, rows 2-43) -- alignment must be set BEFORE the
#        number format below so the date column keeps the built-in
#        short-date format (numFmtId 14) instead of a custom one.
$ws.Range("A2:L43").HorizontalAlignment = -4108

# --- 4. Centered short-date format for the "Date Solved" column on the
#        newly added rows.
$ws.Range("F2:F5").NumberFormat = "mm-dd-yy"

# --- 5. Reset the saved selection back to the top-left cell (the author's
#        saved file no longer highlights the old B11:C18 block).
$ws.Range("A1").Select()
